$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 8 (franzosa_ControlvsDisease_Age), shifting existing rows down
$ws.Rows.Item(8).Insert()

# Insert a new row before what is now row 12 (franzosa_ControlvsUC_Fp), shifting rows down
$ws.Rows.Item(12).Insert()

# Fill new row 8: franzosa_ControlvsCD_ConvCD
$ws.Range("A8").Value = "franzosa_ControlvsCD_ConvCD"
$ws.Range("B8").Value = 0.02
$ws.Range("C8").Value = 0.12
$ws.Range("D8").Value = 0.02
$ws.Range("E8").Value = 0.32
$ws.Range("F8").Value = 0.84
$ws.Range("G8").Value = 0.54
$ws.Range("H8").Value = 0.64

# Fill new row 12: franzosa_ControlvsUC_ConvUC
$ws.Range("A12").Value = "franzosa_ControlvsUC_ConvUC"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0.02
$ws.Range("E12").Value = 0.38
$ws.Range("F12").Value = 0.98
$ws.Range("G12").Value = 0.62
$ws.Range("H12").Value = 0.6
